$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the current column N (14), matching the
# width of the column to its left (M).
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Activate the "Repayment schedule" sheet and select the new cell.
$ws.Activate()
$ws.Range("R8").Select()

Write-Output "done"
